$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16, pushing existing rows 16-29 down to 17-30.
$ws.Rows("16:16").Insert()

# Populate the newly-inserted row 16 with the new weekly price record.
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 45205
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = 100112039
$ws.Range("G16").Value = "Ciboulette"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2250
$ws.Range("N16").Value = "$/docena de atados"
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 750
$ws.Range("Q16").Value = 3
$ws.Range("R16").Value = "Hortaliza"
